# Update ig name
# - Replace the "[code]" placeholder in the IG URLs with "tde"
# - Bump the Version from 0.1.0 to 2.0.0
# - Update the Date timestamp
# - Refresh the matching values on the Elements sheet

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# Metadata sheet: URL, Version, Date
$wsMetadata.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/StructureDefinition/EyeColor"
$wsMetadata.Range("B3").Value = "2.0.0"
$wsMetadata.Range("B8").Value = "2026-01-15T15:23:39+00:00"

# Elements sheet: Fixed Value (URL) and Binding Value Set (ValueSet URL)
$wsElements.Range("R5").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/StructureDefinition/EyeColor"
$wsElements.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/EyeColorVS"

# Re-fit column Z ("Binding Value Set") width to its new (shorter) best-fit
# content width now that the ValueSet URL is a few characters shorter.
$wsElements.Range("Z1").EntireColumn.ColumnWidth = 48.65
